# Update "Predicted Eg" values (column B) on the active worksheet to
# reflect the results of the newly-added Random Forest algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.09
$ws.Range("B3").Value = 3.18
$ws.Range("B4").Value = 3.25
$ws.Range("B5").Value = 3.2
$ws.Range("B6").Value = 3.19
$ws.Range("B7").Value = 3.21
$ws.Range("B8").Value = 3.2
$ws.Range("B9").Value = 3.01
$ws.Range("B10").Value = 3.06
$ws.Range("B11").Value = 3.12
$ws.Range("B12").Value = 3.1
$ws.Range("B13").Value = 3.16
$ws.Range("B14").Value = 3.12
$ws.Range("B16").Value = 3.02
$ws.Range("B17").Value = 3.05
$ws.Range("B18").Value = 3.06
$ws.Range("B19").Value = 3.05
$ws.Range("B20").Value = 1.95
$ws.Range("B21").Value = 2.69
$ws.Range("B22").Value = 2.88
$ws.Range("B23").Value = 2.86
$ws.Range("B24").Value = 2.87
$ws.Range("B25").Value = 2.79
$ws.Range("B26").Value = 5.62
$ws.Range("B27").Value = 3.39
$ws.Range("B28").Value = 3.7
$ws.Range("B29").Value = 2.35
$ws.Range("B30").Value = 3.58
$ws.Range("B31").Value = 8.52
$ws.Range("B32").Value = 2.54
$ws.Range("B33").Value = 2.55
$ws.Range("B34").Value = 3
$ws.Range("B35").Value = 3.19
